$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels: "_old" -> "_FV2304", "_new" -> "_FV2310"
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = [string]$cell.Value2
    if ($val.EndsWith("_old")) {
        $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2304"
    } elseif ($val.EndsWith("_new")) {
        $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2310"
    }
}

# Freeze the header row (row 1)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Convert the data range into an Excel Table
$range = $ws.Range("A1:U72")
$tbl = $ws.ListObjects.Add(1, $range, 0, 1)
$tbl.Name = "Table1"
